# Update cryptocurrency price/volume figures to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each target cell already holds plain text (inline string) in the source
# file. Forcing NumberFormat to "@" (text) before the assignment stops the
# COM layer from re-interpreting number-looking text (e.g. "1.00", "6.11")
# as a numeric value, and ClearFormats() afterwards drops the temporary
# style index again so the cell ends up with no "s" attribute at all -
# exactly like the original file.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '56.956.92'
Set-TextValue $ws.Range('E2') '  -0.77%  '
Set-TextValue $ws.Range('D3') '2.317.32'
Set-TextValue $ws.Range('E3') '  -2.09%  '
Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  +0.28%  '
Set-TextValue $ws.Range('D5') '529.81'
Set-TextValue $ws.Range('E5') '  +2.19%  '
Set-TextValue $ws.Range('D6') '132.44'
Set-TextValue $ws.Range('E6') '  -2.37%  '
Set-TextValue $ws.Range('D7') '0.995'
Set-TextValue $ws.Range('E7') '  -0.11%  '
Set-TextValue $ws.Range('E8') '  -1.16%  '
Set-TextValue $ws.Range('D9') '2.342.00'
Set-TextValue $ws.Range('E9') '  -2.06%  '
Set-TextValue $ws.Range('E10') '  -1.41%  '
Set-TextValue $ws.Range('E11') '  +0.12%  '
Set-TextValue $ws.Range('E12') '  -2.83%  '
Set-TextValue $ws.Range('E13') '  +2.59%  '
Set-TextValue $ws.Range('D14') '2.735.91'
Set-TextValue $ws.Range('E14') '  -2.04%  '
Set-TextValue $ws.Range('D15') '23.32'
Set-TextValue $ws.Range('E15') '  -4.23%  '
Set-TextValue $ws.Range('D16') '56.994.89'
Set-TextValue $ws.Range('E16') '  -0.71%  '
Set-TextValue $ws.Range('E17') '  -2.22%  '
Set-TextValue $ws.Range('D18') '2.341.96'
Set-TextValue $ws.Range('E18') '  -1.51%  '
Set-TextValue $ws.Range('D19') '337.07'
Set-TextValue $ws.Range('E19') '  +2.11%  '
Set-TextValue $ws.Range('D20') '10.40'
Set-TextValue $ws.Range('E20') '  -1.76%  '
Set-TextValue $ws.Range('E21') '  -1.74%  '
Set-TextValue $ws.Range('D22') '6.77'
Set-TextValue $ws.Range('E22') '  +0.86%  '
Set-TextValue $ws.Range('E23') '  -0.11%  '
Set-TextValue $ws.Range('D24') '61.85'
Set-TextValue $ws.Range('E24') '  +0.61%  '
Set-TextValue $ws.Range('D25') '0.167'
Set-TextValue $ws.Range('E25') '  +0.40%  '
Set-TextValue $ws.Range('D26') '8.69'
Set-TextValue $ws.Range('E26') '  -2.89%  '
Set-TextValue $ws.Range('D27') '0.995'
Set-TextValue $ws.Range('E27') '  -0.13%  '
Set-TextValue $ws.Range('E28') '  +1.14%  '
Set-TextValue $ws.Range('D29') '173.34'
Set-TextValue $ws.Range('E30') '  +1.08%  '
Set-TextValue $ws.Range('D31') '0.0₃0724'
Set-TextValue $ws.Range('E31') '  -3.18%  '
Set-TextValue $ws.Range('D32') '6.11'
Set-TextValue $ws.Range('E32') '  -2.55%  '
Set-TextValue $ws.Range('E33') '  -0.69%  '
Set-TextValue $ws.Range('D34') '0.999'
Set-TextValue $ws.Range('E34') '  -0.04%  '
Set-TextValue $ws.Range('D35') '0.992'
Set-TextValue $ws.Range('E35') '  -0.19%  '
Set-TextValue $ws.Range('E36') '  -4.09%  '
Set-TextValue $ws.Range('E37') '  +1.23%  '
Set-TextValue $ws.Range('D38') '3.98'
Set-TextValue $ws.Range('E38') '  -1.68%  '
Set-TextValue $ws.Range('D39') '39.22'
Set-TextValue $ws.Range('E39') '  +0.83%  '
Set-TextValue $ws.Range('D40') '1.58'
Set-TextValue $ws.Range('E40') '  -2.42%  '
Set-TextValue $ws.Range('D41') '5.82'
Set-TextValue $ws.Range('E41') '  +10.04%  '
Set-TextValue $ws.Range('D42') '149.33'
Set-TextValue $ws.Range('E42') '  -0.77%  '
Set-TextValue $ws.Range('D43') '0.375'
Set-TextValue $ws.Range('E43') '  -3.38%  '
Set-TextValue $ws.Range('E44') '  -1.60%  '
Set-TextValue $ws.Range('D45') '283.80'
Set-TextValue $ws.Range('E45') '  -1.90%  '
Set-TextValue $ws.Range('E46') '  -1.20%  '
Set-TextValue $ws.Range('D48') '18.73'
Set-TextValue $ws.Range('E48') '  +2.94%  '
Set-TextValue $ws.Range('E49') '  -1.77%  '
Set-TextValue $ws.Range('E50') '  -1.20%  '
Set-TextValue $ws.Range('E51') '  -2.63%  '
